$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove row 520 ("positive psychotic disorder symptom" / GMHO:0000077),
# shifting all subsequent rows up by one. This reduces the used range
# from A1:K614 to A1:K613.
$ws.Rows.Item(520).Delete()
